$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 290. This pushes the existing rows 290-399
# down to 291-400, preserving their data/formatting.
$ws.Rows.Item(290).Insert()

# Fill the new row 290 with the new weekly record.
$ws.Range("A290").Value = 3
$ws.Range("B290").Value = "Femacal de La Calera"
$ws.Range("C290").Value = "Coquimbo"
$ws.Range("D290").Value = 44755
$ws.Range("D290").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E290").Value = 5
$ws.Range("F290").Value = 100112031
$ws.Range("G290").Value = "Poroto verde"
$ws.Range("H290").Value = "Magnum"
$ws.Range("I290").Value = "Primera"
$ws.Range("J290").Value = 38
$ws.Range("K290").Value = 31000
$ws.Range("L290").Value = 31000
$ws.Range("M290").Value = 31000
$ws.Range("N290").Value = "$/malla 25 kilos"
$ws.Range("O290").Value = "Región de Arica y Parinacota"
$ws.Range("P290").Value = 1240
$ws.Range("Q290").Value = 25
$ws.Range("R290").Value = "Hortaliza"
